# Refresh cached market-board profit figures on the per-job Leve tables.
# Source data refreshed by the scheduled runner; only numeric market/profit
# columns (H..N: currentAveragePrice*, LevePrice*, LeveProfit*) are touched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 3856.9443
$ws.Range("J62").Value = 4574.6665
$ws.Range("L62").Value = 4574.6665
$ws.Range("N62").Value = -5822.6665

# Row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 3856.9443
$ws.Range("J65").Value = 4574.6665
$ws.Range("L65").Value = 22873.3325
$ws.Range("N65").Value = -29113.3325

# Row 107: Another Man's Ink / Enchanted Truegold Ink
$ws.Range("H107").Value = 634.7619
$ws.Range("I107").Value = 606.8889
$ws.Range("K107").Value = 606.8889
$ws.Range("M107").Value = 1313.1111

# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 3312.2163
$ws.Range("J112").Value = 3373.6667
$ws.Range("L112").Value = 10121.0001
$ws.Range("N112").Value = -12337.0001

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 1984.2
$ws.Range("I32").Value = 1797.1632
$ws.Range("K32").Value = 1797.1632
$ws.Range("M32").Value = -1510.1632

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 5712.7856
$ws.Range("I122").Value = 3275.4443
$ws.Range("K122").Value = 9826.332900000001
$ws.Range("M122").Value = -7376.332900000001

$ws = $wb.Worksheets.Item("BSM")
# Row 42: Hard Knock Life / Steel Sledgehammer
$ws.Range("H42").Value = 232497.67
$ws.Range("J42").Value = 232497.67
$ws.Range("L42").Value = 232497.67
$ws.Range("N42").Value = -233153.67

$ws = $wb.Worksheets.Item("CRP")
# Row 12: A Sword in Hand / Ash Macuahuitl
$ws.Range("H12").Value = 419
$ws.Range("I12").Value = 419
$ws.Range("K12").Value = 419
$ws.Range("M12").Value = -249

# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 28089.191
$ws.Range("I62").Value = 21128.076
$ws.Range("J62").Value = 39401
$ws.Range("K62").Value = 21128.076
$ws.Range("L62").Value = 39401
$ws.Range("M62").Value = -20504.076
$ws.Range("N62").Value = -40649

# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 28089.191
$ws.Range("I65").Value = 21128.076
$ws.Range("J65").Value = 39401
$ws.Range("K65").Value = 105640.38
$ws.Range("L65").Value = 197005
$ws.Range("M65").Value = -102520.38
$ws.Range("N65").Value = -203245

$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up / Kukuru Butter
$ws.Range("H12").Value = 144.5
$ws.Range("I12").Value = 227
$ws.Range("J12").Value = 128
$ws.Range("K12").Value = 681
$ws.Range("L12").Value = 384
$ws.Range("M12").Value = -508
$ws.Range("N12").Value = -730

# Row 23: Sweet Smell of Success / Lavender Oil
$ws.Range("H23").Value = 472
$ws.Range("J23").Value = 491.1111
$ws.Range("L23").Value = 1473.3333
$ws.Range("N23").Value = -1943.3333

# Row 42: Point Them with the Sticky End / Tuna Miq'abob
$ws.Range("H42").Value = 14949.5
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 14949.5
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 44848.5
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -45916.5

# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 2079
$ws.Range("I68").Value = 934.6
$ws.Range("J68").Value = 2896.4285
$ws.Range("K68").Value = 2803.8
$ws.Range("L68").Value = 8689.2855
$ws.Range("M68").Value = -1992.8
$ws.Range("N68").Value = -10311.2855

# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 2079
$ws.Range("I71").Value = 934.6
$ws.Range("J71").Value = 2896.4285
$ws.Range("K71").Value = 8411.4
$ws.Range("L71").Value = 26067.8565
$ws.Range("M71").Value = -4355.4
$ws.Range("N71").Value = -34179.8565

# Row 107: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 952.86664
$ws.Range("J107").Value = 861.75
$ws.Range("L107").Value = 2585.25
$ws.Range("N107").Value = -6425.25

# Row 110: His Dark Utensils / Spaghetti al Nero
$ws.Range("H110").Value = 9363.777
$ws.Range("I110").Value = 8455.2
$ws.Range("J110").Value = 10499.5
$ws.Range("K110").Value = 25365.6
$ws.Range("L110").Value = 31498.5
$ws.Range("M110").Value = -21275.6
$ws.Range("N110").Value = -39678.5

# Row 126: Imperial Palate / Glory Be Soup
$ws.Range("H126").Value = 16207.5
$ws.Range("I126").Value = 2418
$ws.Range("J126").Value = 29997
$ws.Range("K126").Value = 7254
$ws.Range("L126").Value = 89991
$ws.Range("M126").Value = -2314
$ws.Range("N126").Value = -99871

# Row 138: Bring Me Your Tacos / Tacos Al Pastor
$ws.Range("H138").Value = 2782.087
$ws.Range("I138").Value = 2759.6
$ws.Range("J138").Value = 2932
$ws.Range("K138").Value = 8278.8
$ws.Range("L138").Value = 8796
$ws.Range("M138").Value = -3138.799999999999
$ws.Range("N138").Value = -19076

# Row 139: Najoothie / Wild Banana Blend
$ws.Range("H139").Value = 2515.4517
$ws.Range("I139").Value = 1974.1666
$ws.Range("J139").Value = 4371.2856
$ws.Range("K139").Value = 5922.4998
$ws.Range("L139").Value = 13113.8568
$ws.Range("M139").Value = -782.4997999999996
$ws.Range("N139").Value = -23393.8568

$ws = $wb.Worksheets.Item("GSM")
# Row 3: Needful Rings / Copper Wristlets
$ws.Range("H3").Value = 23333500
$ws.Range("I3").Value = 35000000
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 35000000
$ws.Range("L3").Value = 500
$ws.Range("M3").Value = -34999884
$ws.Range("N3").Value = -732

# Row 14: All That Glitters / Copper Ear Cuffs
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

# Row 18: Gorgeous Gorget / Brass Gorget
$ws.Range("H18").Value = 18000
$ws.Range("I18").Value = 18000
$ws.Range("K18").Value = 18000
$ws.Range("M18").Value = -17707

# Row 19: Better Four Eyes than None / Brass Spectacles
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

# Row 21: Forever 21K / Brass Ring
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

# Row 30: Dog Tags Are for Dogs / Brass Ring
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 4371.2085
$ws.Range("I7").Value = 4999.25
$ws.Range("K7").Value = 4999.25
$ws.Range("M7").Value = -4887.25

# Row 23: Back in the Band / Hard Leather Wristbands
$ws.Range("H23").Value = 7999.75
$ws.Range("I23").Value = 6999.5
$ws.Range("J23").Value = 9000
$ws.Range("K23").Value = 6999.5
$ws.Range("L23").Value = 9000
$ws.Range("M23").Value = -6769.5
$ws.Range("N23").Value = -9460

# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 3373
$ws.Range("J93").Value = 4104
$ws.Range("L93").Value = 4104
$ws.Range("N93").Value = -6600

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 4371.2085
$ws.Range("I126").Value = 4999.25
$ws.Range("K126").Value = 14997.75
$ws.Range("M126").Value = -12527.75

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 4885.9185
$ws.Range("I132").Value = 3950.6516
$ws.Range("K132").Value = 11851.9548
$ws.Range("M132").Value = -9321.9548

$ws = $wb.Worksheets.Item("WVR")
# Row 17: Making Gloves Out of Nothing at All / Hempen Bracers
$ws.Range("H17").Value = 5600
$ws.Range("J17").Value = 1200
$ws.Range("L17").Value = 1200
$ws.Range("N17").Value = -1544

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 6174343
$ws.Range("I132").Value = 10101628
$ws.Range("K132").Value = 30304884
$ws.Range("M132").Value = -30302354
